# Update cryptos list (price / volume(1h) columns, and two row re-orderings)
# to reflect the refreshed data snapshot from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Subscript-three character used inside one of the PEPE price strings (0.0\u20830778).
$sub3 = [char]0x2083

$ws.Range("D2").Value = "'60.469.84"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "'2.639.57"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'570.51"
$ws.Range("E5").Value = "  +6.57%  "
$ws.Range("D6").Value = "'146.72"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "'0.609"
$ws.Range("E8").Value = "  +7.50%  "
$ws.Range("D9").Value = "'6.83"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("E11").Value = "  +6.35%  "
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "'3.111.49"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'60.459.74"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").Value = "'21.75"
$ws.Range("E15").Value = "  +4.82%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000137"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.652.33"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  +3.78%  "
$ws.Range("D19").Value = "'345.02"
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("D20").Value = "'10.44"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").Value = "'6.39"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").Value = "'5.82"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'66.87"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "'0.444"
$ws.Range("E25").Value = "  +6.89%  "
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "'7.35"
$ws.Range("E28").Value = "  +3.60%  "
$ws.Range("D29").Value = "'0.0" + $sub3 + "0778"
$ws.Range("E29").Value = "  +5.72%  "
$ws.Range("D31").Value = "'1.72"
$ws.Range("E31").Value = "  +4.66%  "
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("D33").Value = "'156.09"
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").Value = "'19.22"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("E35").Value = "  +5.37%  "
$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "'0.912"
$ws.Range("E36").Value = "  +7.29%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'0.911"
$ws.Range("E37").Value = "  +12.62%  "
$ws.Range("E38").Value = "  +6.48%  "
$ws.Range("D39").Value = "'37.62"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  +7.86%  "
$ws.Range("D41").Value = "'306.86"
$ws.Range("E41").Value = "  +9.53%  "
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("D43").Value = "'0.994"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "'0.609"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "'0.0980"
$ws.Range("E45").Value = "  +4.77%  "
$ws.Range("D46").Value = "'0.0549"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("D47").Value = "'19.47"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("D50").Value = "'125.57"
$ws.Range("E50").Value = "  +11.51%  "
$ws.Range("D51").Value = "'1.972.14"
$ws.Range("E51").Value = "  +1.37%  "
